$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$t.Cell(1, 1).Range.Text = "63 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "6|    |" + $vt + "3|    |"
$t.Cell(1, 2).Range.Text = "46 x 64" + $vt + "  6    4" + $vt + "  ----" + $vt + "4|    |" + $vt + "6|    |"
$t.Cell(1, 3).Range.Text = "62 x 57" + $vt + "  5    7" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"
$t.Cell(2, 1).Range.Text = "11 x 92" + $vt + "  9    2" + $vt + "  ----" + $vt + "1|    |" + $vt + "1|    |"
$t.Cell(2, 2).Range.Text = "91 x 30" + $vt + "  3    0" + $vt + "  ----" + $vt + "9|    |" + $vt + "1|    |"
$t.Cell(2, 3).Range.Text = "45 x 84" + $vt + "  8    4" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"
$t.Cell(3, 1).Range.Text = "77 x 97" + $vt + "  9    7" + $vt + "  ----" + $vt + "7|    |" + $vt + "7|    |"
$t.Cell(3, 2).Range.Text = "60 x 10" + $vt + "  1    0" + $vt + "  ----" + $vt + "6|    |" + $vt + "0|    |"
$t.Cell(3, 3).Range.Text = "51 x 15" + $vt + "  1    5" + $vt + "  ----" + $vt + "5|    |" + $vt + "1|    |"
$t.Cell(4, 1).Range.Text = "91 x 26" + $vt + "  2    6" + $vt + "  ----" + $vt + "9|    |" + $vt + "1|    |"
$t.Cell(4, 2).Range.Text = "45 x 55" + $vt + "  5    5" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"
$t.Cell(4, 3).Range.Text = "99 x 58" + $vt + "  5    8" + $vt + "  ----" + $vt + "9|    |" + $vt + "9|    |"
$t.Cell(5, 1).Range.Text = "27 x 46" + $vt + "  4    6" + $vt + "  ----" + $vt + "2|    |" + $vt + "7|    |"
$t.Cell(5, 2).Range.Text = "50 x 39" + $vt + "  3    9" + $vt + "  ----" + $vt + "5|    |" + $vt + "0|    |"
$t.Cell(5, 3).Range.Text = "37 x 67" + $vt + "  6    7" + $vt + "  ----" + $vt + "3|    |" + $vt + "7|    |"
Write-Output "done"
